$wb = $excel.ActiveWorkbook

# Add the new "CardInformation" worksheet right after the existing
# TestDataSheet1 sheet and make it the active sheet (tab).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CardInformation"

# Header row
$ws.Range("A1").Value = "CardType"
$ws.Range("B1").Value = "CardNumber"
$ws.Range("C1").Value = "CVV"
$ws.Range("D1").Value = "ExpDate"

# Data rows - card test data for cart/checkout test cases
$ws.Range("A2").Value = "GiftCard"
$ws.Range("B2").Value = 98990988877
$ws.Range("C2").Value = 202
$ws.Range("D2").Value = 46266

$ws.Range("A3").Value = "CreditCard"
$ws.Range("B3").Value = 58455544115
$ws.Range("C3").Value = 111
$ws.Range("D3").Value = 47481

# ExpDate column is formatted as a month/year date
$ws.Range("D1:D3").NumberFormat = "mmm-yy"

# Approximate the column layout used for this sheet
$ws.Columns.Item(1).ColumnWidth = 11.14
$ws.Columns.Item(2).ColumnWidth = 12.86

# Match the page margins used on the new sheet (inches -> points)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
$ws.PageSetup.LeftHeader = ""

# Selection state left on the new sheet
$null = $ws.Range("D8").Select()
